$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 212
$ws.Cells.Item(3, 6).Value = 10
$ws.Cells.Item(4, 6).Value = 398
$ws.Cells.Item(6, 6).Value = 782
$ws.Cells.Item(7, 6).Value = 84
$ws.Cells.Item(8, 6).Value = 10075
$ws.Cells.Item(9, 6).Value = 56
$ws.Cells.Item(10, 6).Value = 3467
$ws.Cells.Item(12, 6).Value = 2427
$ws.Cells.Item(14, 6).Value = 2766
$ws.Cells.Item(16, 6).Value = 506
$ws.Cells.Item(17, 6).Value = 2140
$ws.Cells.Item(19, 6).Value = 92
$ws.Cells.Item(20, 6).Value = 18
$ws.Cells.Item(23, 6).Value = 128
$ws.Cells.Item(25, 6).Value = 271
$ws.Cells.Item(26, 6).Value = 204
$ws.Cells.Item(27, 6).Value = 612
$ws.Cells.Item(28, 6).Value = 1308
$ws.Cells.Item(29, 6).Value = 4
$ws.Cells.Item(30, 6).Value = 1247
$ws.Cells.Item(34, 6).Value = 2731
$ws.Cells.Item(35, 6).Value = 2949
$ws.Cells.Item(36, 6).Value = 17
$ws.Cells.Item(38, 6).Value = 1025
$ws.Cells.Item(39, 6).Value = 386
$ws.Cells.Item(41, 6).Value = 1288
$ws.Cells.Item(42, 6).Value = 85
$ws.Cells.Item(43, 6).Value = 102
$ws.Cells.Item(44, 6).Value = 67
$ws.Cells.Item(47, 6).Value = 6

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 174
$ws.Cells.Item(8, 6).Value = 5

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 739
$ws.Cells.Item(3, 6).Value = 977
$ws.Cells.Item(4, 6).Value = 124
$ws.Cells.Item(5, 6).Value = 1972

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 739
$ws.Cells.Item(3, 6).Value = 977
$ws.Cells.Item(4, 6).Value = 10
$ws.Cells.Item(5, 6).Value = 398
$ws.Cells.Item(9, 6).Value = 782
$ws.Cells.Item(10, 6).Value = 84
$ws.Cells.Item(11, 6).Value = 10075
$ws.Cells.Item(12, 6).Value = 174
$ws.Cells.Item(13, 6).Value = 56
$ws.Cells.Item(15, 6).Value = 3468
$ws.Cells.Item(16, 6).Value = 2427
$ws.Cells.Item(18, 6).Value = 2766
$ws.Cells.Item(20, 6).Value = 506
$ws.Cells.Item(21, 6).Value = 2140
$ws.Cells.Item(23, 6).Value = 92
$ws.Cells.Item(24, 6).Value = 128
$ws.Cells.Item(26, 6).Value = 271
$ws.Cells.Item(27, 6).Value = 612
$ws.Cells.Item(28, 6).Value = 1308
$ws.Cells.Item(29, 6).Value = 1247
$ws.Cells.Item(31, 6).Value = 5
$ws.Cells.Item(33, 6).Value = 2731
$ws.Cells.Item(35, 6).Value = 2949
$ws.Cells.Item(36, 6).Value = 1025
$ws.Cells.Item(39, 6).Value = 386
$ws.Cells.Item(44, 6).Value = 1288
$ws.Cells.Item(45, 6).Value = 85
$ws.Cells.Item(46, 6).Value = 67
